# daily auto push: 2026-02-17 03:13 UTC
# Insert a new data row (835) into the "sei1" time-series sheet, pushing the
# existing rows 835..876 down to 836..877.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 835..876 down to 836..877 and leave row 835 blank.
$ws.Rows.Item(835).Insert()

# Column A ("日付" / date) holds plain text like "2026/02/17", never a real
# Excel date value (that's how every other row in the sheet stores it).
# Assigning a date-looking string directly would make Excel auto-convert it
# to a date serial number, so force the cell to Text format first, write the
# string, then drop the number-format override again so the cell ends up
# with no special styling - exactly like its neighbours.
$ws.Cells.Item(835, 1).NumberFormat = "@"
$ws.Cells.Item(835, 1).Value = "2026/02/17"
$ws.Cells.Item(835, 1).ClearFormats()

# Column B ("曜日" / weekday) - plain text, no special handling needed.
$ws.Cells.Item(835, 2).Value = "火"

# Columns C ("時刻") and D ("ランキング") are plain numbers.
$ws.Cells.Item(835, 3).Value = 10
$ws.Cells.Item(835, 4).Value = 201
